$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 5, shifting existing rows 5-23 down to 6-24
$ws.Rows.Item(5).Insert()

# Populate the newly inserted row with the new leaderboard entry
$ws.Cells.Item(5, 1).Value = "aulop"
$ws.Cells.Item(5, 2).Value = 2
$ws.Cells.Item(5, 3).Value = "2025-12-18 22:13:30"
$ws.Cells.Item(5, 4).Value = "Normal"
